$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Plan")

# Game Over section (rows 27-30 & 34): Sprint 3 status flipped Fail -> Pass,
# and tester comments amended with SJD's follow-up note.
$ws.Range("G27").Value = "Pass"
$ws.Range("H27").Value = "AO; 4/5/18; Winning doesn't work .Got rid of all white pieces and nothing happened; SJD; 4/7/18"

$ws.Range("G29").Value = "Pass"
$ws.Range("H29").Value = "AO; 4/5/18; see above; SJD; 4/7/18"

$ws.Range("G30").Value = "Pass"
$ws.Range("H30").Value = "AO; 4/5/18; see above; SJD; 4/7/18"

$ws.Range("G34").Value = "Pass"
$ws.Range("H34").Value = "AO; 4/5/18; If it is the Player's turn, they are not redirected to the home page until they submit their turn (change your acceptance criteria); SJD 4/7/18; It must operate that way to stay within the vision document"

# New Player Help acceptance criteria rows (47-49).
$ws.Range("B47").Value = "Given I am a player when I am in a game then I expect to have  a link to request help."
$ws.Range("B48").Value = "Given I am a player when I click the request help link then I expect to be redirected to a help page."
$ws.Range("B49").Value = "Given I am a player when I am on the help page then I expect to be able to return to my game."
